# Apply the "feat: add 2022-Q1 data" change:
#  - The existing "总计" (Total) sheet becomes the new "2022-Q1" sheet
#    (its old total-by-quarter content is replaced with the new quarter's
#    per-fund holdings).
#  - A brand-new "总计" sheet is appended at the end with the refreshed
#    totals table (the new 2022-Q1 row plus the previous rows).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Turn the current "总计" sheet into "2022-Q1"
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Use an existing per-fund sheet as a formatting template for the header
# row and the index column.
$template = $wb.Worksheets.Item("2021-Q4")

$q1.Cells.Clear()

# Copy header-row formatting (B1:H1) and index-column formatting (A2:A3)
$template.Range("B1:H1").Copy() | Out-Null
$q1.Range("B1:H1").PasteSpecial(-4122) | Out-Null
$template.Range("A2:A3").Copy() | Out-Null
$q1.Range("A2:A3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$headers = New-Object 'object[,]' 1,7
$headers[0,0] = "基金代码"
$headers[0,1] = "基金名称"
$headers[0,2] = "基金规模"
$headers[0,3] = "股票总仓位"
$headers[0,4] = "仓位占比"
$headers[0,5] = "持有市值(亿元)"
$headers[0,6] = "仓位排名"
$q1.Range("B1:H1").Value = $headers

# Numeric-looking values (fund code / fund size / position figures) must
# stay text, matching how this source data is stored everywhere else in
# the workbook - a leading apostrophe forces text entry, just like a user
# typing into the formula bar would.
$row2 = New-Object 'object[,]' 1,7
$row2[0,0] = "'004854"
$row2[0,1] = "广发中证全指汽车指数A"
$row2[0,2] = "'22.01"
$row2[0,3] = "'94.43"
$row2[0,4] = "'2.72"
$row2[0,5] = "'0.5987"
$row2[0,6] = 9
$q1.Range("A2").Value = 0
$q1.Range("B2:H2").Value = $row2

$row3 = New-Object 'object[,]' 1,7
$row3[0,0] = "'004855"
$row3[0,1] = "广发中证全指汽车指数C"
$row3[0,2] = "'6.11"
$row3[0,3] = "'94.43"
$row3[0,4] = "'2.72"
$row3[0,5] = "'0.1662"
$row3[0,6] = 9
$q1.Range("A3").Value = 1
$q1.Range("B3:H3").Value = $row3

# ---------------------------------------------------------------------
# 2. Append a fresh "总计" sheet after "2022-Q1" with the updated totals
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1)
$total.Name = "总计"

$totalTemplate = $wb.Worksheets.Item("2021-Q4")
$totalTemplate.Range("B1:D1").Copy() | Out-Null
$total.Range("B1:D1").PasteSpecial(-4122) | Out-Null
$totalTemplate.Range("A2:A3").Copy() | Out-Null
$total.Range("A2:A7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$totalHeaders = New-Object 'object[,]' 1,3
$totalHeaders[0,0] = "日期"
$totalHeaders[0,1] = "持有数量(只)"
$totalHeaders[0,2] = "持有市值(亿元)"
$total.Range("B1:D1").Value = $totalHeaders

$data = @(
    @(0, "2022-Q1", 2, 0.76),
    @(1, "2021-Q4", 3, 0.2),
    @(2, "2021-Q3", 8, 4.47),
    @(3, "2021-Q2", 8, 3.53),
    @(4, "2021-Q1", 10, 3.82),
    @(5, "2020-Q4", 6, 4.68)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    $total.Range("A$r").Value = $rowVals[0]
    $rowData = New-Object 'object[,]' 1,3
    $rowData[0,0] = $rowVals[1]
    $rowData[0,1] = $rowVals[2]
    $rowData[0,2] = $rowVals[3]
    $total.Range("B$r`:D$r").Value = $rowData
}

Write-Host "done"
